# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" / "Valor Mora" rows (16-22) are re-sorted into
# ascending period order. Previously the periods ran 1908 .. 1901
# (descending); now they run 1901 .. 1908 (ascending), carrying their
# "Valor Mora" along with them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ascending period order for rows 16..22.
$periodos = @("1901", "1903", "1904", "1905", "1906", "1907", "1908")
$valores  = @(31249, 33125, 33125, 33125, 33125, 33125, 20979)

for ($i = 0; $i -lt 7; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periodos[$i]
    $ws.Range("F$row").Value = $valores[$i]
}
